$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 969.55
$ws.Range("I15").Value = 969.55
$ws.Range("K15").Value = 2908.65
$ws.Range("M15").Value = -2739.65

$ws.Range("H40").Value = 1908.8148
$ws.Range("I40").Value = 1757.8462
$ws.Range("J40").Value = 2049
$ws.Range("K40").Value = 1757.8462
$ws.Range("L40").Value = 2049
$ws.Range("M40").Value = -1582.8462
$ws.Range("N40").Value = -2399

$ws.Range("H64").Value = 3819.9167
$ws.Range("J64").Value = 3993.3333
$ws.Range("L64").Value = 3993.3333
$ws.Range("N64").Value = -4489.3333

$ws.Range("H67").Value = 3819.9167
$ws.Range("J67").Value = 3993.3333
$ws.Range("L67").Value = 3993.3333
$ws.Range("N67").Value = -5709.3333

$ws.Range("H70").Value = 1854.2858
$ws.Range("I70").Value = 1808.5714
$ws.Range("J70").Value = 1900
$ws.Range("K70").Value = 5425.7142
$ws.Range("L70").Value = 5700
$ws.Range("M70").Value = -5155.7142
$ws.Range("N70").Value = -6240

$ws.Range("H73").Value = 1854.2858
$ws.Range("I73").Value = 1808.5714
$ws.Range("J73").Value = 1900
$ws.Range("K73").Value = 5425.7142
$ws.Range("L73").Value = 5700
$ws.Range("M73").Value = -4489.7142
$ws.Range("N73").Value = -7572

$ws.Range("H74").Value = 3250
$ws.Range("I74").Value = 2500
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 2500
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -1564
$ws.Range("N74").Value = -5872

$ws.Range("H76").Value = 6620
$ws.Range("I76").Value = 6620
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 6620
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -6305

$ws.Range("H77").Value = 3250
$ws.Range("I77").Value = 2500
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 12500
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -7820
$ws.Range("N77").Value = -29360

$ws.Range("H79").Value = 6620
$ws.Range("I79").Value = 6620
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 6620
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -5528

$ws.Range("H100").Value = 1038.6364
$ws.Range("I100").Value = 688.8889
$ws.Range("K100").Value = 688.8889
$ws.Range("M100").Value = -147.8889

$ws.Range("H116").Value = 2596.5
$ws.Range("I116").Value = 1753
$ws.Range("J116").Value = 3861.75
$ws.Range("K116").Value = 1753
$ws.Range("L116").Value = 3861.75
$ws.Range("M116").Value = 1689
$ws.Range("N116").Value = -10745.75

$ws.Range("H137").Value = 1174.2963
$ws.Range("I137").Value = 993.3103599999999
$ws.Range("J137").Value = 1384.24
$ws.Range("K137").Value = 2979.93108
$ws.Range("L137").Value = 4152.72
$ws.Range("M137").Value = -429.9310799999998
$ws.Range("N137").Value = -9252.720000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 698.8
$ws.Range("I4").Value = 390
$ws.Range("K4").Value = 390
$ws.Range("M4").Value = -274

$ws.Range("H32").Value = 2400.81
$ws.Range("I32").Value = 1998.6967
$ws.Range("K32").Value = 1998.6967
$ws.Range("M32").Value = -1711.6967

$ws.Range("H61").Value = 76924184
$ws.Range("I61").Value = 83334290
$ws.Range("K61").Value = 83334290
$ws.Range("M61").Value = -83334078

$ws.Range("H74").Value = 1754.5625
$ws.Range("I74").Value = 1359.6154
$ws.Range("K74").Value = 1359.6154
$ws.Range("M74").Value = -485.6153999999999

$ws.Range("H77").Value = 1754.5625
$ws.Range("I77").Value = 1359.6154
$ws.Range("K77").Value = 6798.076999999999
$ws.Range("M77").Value = -2430.076999999999

$ws.Range("H132").Value = 1920.2678
$ws.Range("I132").Value = 1539.5834
$ws.Range("J132").Value = 2605.5
$ws.Range("K132").Value = 4618.7502
$ws.Range("L132").Value = 7816.5
$ws.Range("M132").Value = -2088.7502
$ws.Range("N132").Value = -12876.5

$ws.Range("H136").Value = 76924184
$ws.Range("I136").Value = 83334290
$ws.Range("K136").Value = 250002870
$ws.Range("M136").Value = -250000320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1657.7333
$ws.Range("I31").Value = 1540.4474
$ws.Range("J31").Value = 2294.4285
$ws.Range("K31").Value = 1540.4474
$ws.Range("L31").Value = 2294.4285
$ws.Range("M31").Value = -1245.4474
$ws.Range("N31").Value = -2884.4285

$ws.Range("H34").Value = 1657.7333
$ws.Range("I34").Value = 1540.4474
$ws.Range("J34").Value = 2294.4285
$ws.Range("K34").Value = 1540.4474
$ws.Range("L34").Value = 2294.4285
$ws.Range("M34").Value = -1338.4474
$ws.Range("N34").Value = -2698.4285

$ws.Range("H58").Value = 1448.262
$ws.Range("I58").Value = 1082.3572
$ws.Range("J58").Value = 2180.0715
$ws.Range("K58").Value = 1082.3572
$ws.Range("L58").Value = 2180.0715
$ws.Range("M58").Value = -879.3571999999999
$ws.Range("N58").Value = -2586.0715

$ws.Range("H86").Value = 4806751.5
$ws.Range("I86").Value = 8372020.5
$ws.Range("J86").Value = 53059.5
$ws.Range("K86").Value = 8372020.5
$ws.Range("L86").Value = 53059.5
$ws.Range("M86").Value = -8370897.5
$ws.Range("N86").Value = -55305.5

$ws.Range("H89").Value = 4806751.5
$ws.Range("I89").Value = 8372020.5
$ws.Range("J89").Value = 53059.5
$ws.Range("K89").Value = 41860102.5
$ws.Range("L89").Value = 265297.5
$ws.Range("M89").Value = -41854486.5
$ws.Range("N89").Value = -276529.5

$ws.Range("H109").Value = 9450.5
$ws.Range("J109").Value = 9450.5
$ws.Range("L109").Value = 9450.5
$ws.Range("N109").Value = -11530.5

$ws.Range("H136").Value = 1448.262
$ws.Range("I136").Value = 1082.3572
$ws.Range("J136").Value = 2180.0715
$ws.Range("K136").Value = 3247.0716
$ws.Range("L136").Value = 6540.2145
$ws.Range("M136").Value = -697.0715999999998
$ws.Range("N136").Value = -11640.2145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4877.727
$ws.Range("J107").Value = 5547.8945
$ws.Range("L107").Value = 16643.6835
$ws.Range("N107").Value = -20483.6835

$ws.Range("H122").Value = 1714.3334
$ws.Range("I122").Value = 863.5
$ws.Range("J122").Value = 1957.4286
$ws.Range("K122").Value = 7771.5
$ws.Range("L122").Value = 17616.8574
$ws.Range("M122").Value = -5321.5
$ws.Range("N122").Value = -22516.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 34619016
$ws.Range("I70").Value = 35717884
$ws.Range("J70").Value = 33337000
$ws.Range("K70").Value = 35717884
$ws.Range("L70").Value = 33337000
$ws.Range("M70").Value = -35717614
$ws.Range("N70").Value = -33337540

$ws.Range("H73").Value = 34619016
$ws.Range("I73").Value = 35717884
$ws.Range("J73").Value = 33337000
$ws.Range("K73").Value = 35717884
$ws.Range("L73").Value = 33337000
$ws.Range("M73").Value = -35716948
$ws.Range("N73").Value = -33338872

$ws.Range("H102").Value = 1346.6842
$ws.Range("I102").Value = 1192.5
$ws.Range("J102").Value = 1518
$ws.Range("K102").Value = 1192.5
$ws.Range("L102").Value = 1518
$ws.Range("M102").Value = 429.5
$ws.Range("N102").Value = -4762

$ws.Range("H123").Value = 21400
$ws.Range("J123").Value = 21400
$ws.Range("L123").Value = 21400
$ws.Range("N123").Value = -26300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1249
$ws.Range("I61").Value = 1170.2858
$ws.Range("J61").Value = 1800
$ws.Range("K61").Value = 1170.2858
$ws.Range("L61").Value = 1800
$ws.Range("M61").Value = -968.2858000000001
$ws.Range("N61").Value = -2204

$ws.Range("H68").Value = 1885.5758
$ws.Range("I68").Value = 1850.75
$ws.Range("K68").Value = 1850.75
$ws.Range("M68").Value = -1101.75

$ws.Range("H71").Value = 1885.5758
$ws.Range("I71").Value = 1850.75
$ws.Range("K71").Value = 9253.75
$ws.Range("M71").Value = -5509.75

$ws.Range("H82").Value = 1698.091
$ws.Range("I82").Value = 1621.6471
$ws.Range("J82").Value = 1958
$ws.Range("K82").Value = 1621.6471
$ws.Range("L82").Value = 1958
$ws.Range("M82").Value = -1260.6471
$ws.Range("N82").Value = -2680

$ws.Range("H85").Value = 1698.091
$ws.Range("I85").Value = 1621.6471
$ws.Range("J85").Value = 1958
$ws.Range("K85").Value = 1621.6471
$ws.Range("L85").Value = 1958
$ws.Range("M85").Value = -373.6470999999999
$ws.Range("N85").Value = -4454

$ws.Range("H93").Value = 800
$ws.Range("I93").Value = 800
$ws.Range("K93").Value = 800
$ws.Range("M93").Value = 448

$ws.Range("H113").Value = 1249
$ws.Range("I113").Value = 1170.2858
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1170.2858
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 999.7141999999999
$ws.Range("N113").Value = -6140

$ws.Range("H122").Value = 22729572
$ws.Range("I122").Value = 31252300
$ws.Range("K122").Value = 93756900
$ws.Range("M122").Value = -93754450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1319.05
$ws.Range("I96").Value = 1372
$ws.Range("J96").Value = 1160.2
$ws.Range("K96").Value = 1372
$ws.Range("L96").Value = 1160.2
$ws.Range("M96").Value = 1
$ws.Range("N96").Value = -3906.2

$ws.Range("H122").Value = 12501496
$ws.Range("I122").Value = 14707302
$ws.Range("J122").Value = 1926.6666
$ws.Range("K122").Value = 44121906
$ws.Range("L122").Value = 5779.9998
$ws.Range("M122").Value = -44119456
